$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy formatting from column F to the new column G
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F2:F21").Copy()
$ws.Range("G2:G21").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F22").Copy()
$ws.Range("G22").PasteSpecial(-4122)  # xlPasteFormats

# Set the new column width (COM ColumnWidth uses a slightly different unit than
# the raw OOXML "width" attribute; 16.14 round-trips to a stored width of 17)
$ws.Columns.Item(7).ColumnWidth = 16.14

# Header
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# Budget values for rows 2-21 (column G)
$values = @(0, 1000, 0, 3000, 5000, 0, 6000, 0, 7000, 0, 6000, 6000, 1000, 400, 6500, 0, 4000, 0, 500, 4000)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}

# Total row 22
$ws.Cells.Item(22, 7).Value = 50400
